$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix comment timestamp strings in column M (Creation date): 2-digit year -> 4-digit year
$ws.Range("M2").Value = "8/3/2018 19:27:00"
$ws.Range("M3").Value = "8/3/2018 19:28:00"
$ws.Range("M4").Value = "8/3/2018 19:28:00"
$ws.Range("M5").Value = "8/3/2018 19:28:00"
$ws.Range("M6").Value = "8/3/2018 19:29:00"
$ws.Range("M7").Value = "8/3/2018 19:29:00"
$ws.Range("M8").Value = "8/3/2018 19:29:00"
$ws.Range("M9").Value = "8/3/2018 19:31:00"
$ws.Range("M10").Value = "8/3/2018 19:32:00"
$ws.Range("M11").Value = "8/3/2018 19:32:00"
$ws.Range("M12").Value = "8/3/2018 19:32:00"
$ws.Range("M13").Value = "8/3/2018 19:32:00"
$ws.Range("M14").Value = "8/3/2018 19:33:00"
$ws.Range("M15").Value = "8/3/2018 19:33:00"
$ws.Range("M16").Value = "8/3/2018 19:33:00"
$ws.Range("M17").Value = "8/3/2018 19:34:00"
$ws.Range("M18").Value = "8/3/2018 19:34:00"
$ws.Range("M19").Value = "8/3/2018 19:35:00"
$ws.Range("M20").Value = "8/3/2018 19:35:00"
$ws.Range("M21").Value = "8/3/2018 19:35:00"
$ws.Range("M22").Value = "8/3/2018 19:35:00"
$ws.Range("M23").Value = "8/3/2018 19:36:00"
$ws.Range("M24").Value = "8/3/2018 19:36:00"
$ws.Range("M25").Value = "8/3/2018 19:36:00"
$ws.Range("M26").Value = "8/3/2018 19:37:00"
$ws.Range("M27").Value = "8/3/2018 19:37:00"
$ws.Range("M28").Value = "8/3/2018 19:38:00"
$ws.Range("M29").Value = "8/3/2018 19:38:00"
$ws.Range("M30").Value = "8/3/2018 19:39:00"
$ws.Range("M31").Value = "8/3/2018 19:41:00"
$ws.Range("M32").Value = "8/3/2018 19:41:00"
$ws.Range("M33").Value = "8/3/2018 19:41:00"
$ws.Range("M34").Value = "8/3/2018 19:42:00"
$ws.Range("M35").Value = "8/3/2018 19:42:00"
$ws.Range("M36").Value = "8/3/2018 19:43:00"
$ws.Range("M37").Value = "8/3/2018 19:43:00"
$ws.Range("M38").Value = "8/3/2018 19:43:00"
$ws.Range("M39").Value = "8/3/2018 19:43:00"
$ws.Range("M40").Value = "8/3/2018 19:43:00"
$ws.Range("M41").Value = "8/3/2018 19:43:00"
$ws.Range("M42").Value = "8/3/2018 19:44:00"
$ws.Range("M43").Value = "8/3/2018 19:44:00"
$ws.Range("M44").Value = "8/3/2018 19:45:00"
$ws.Range("M45").Value = "8/3/2018 19:45:00"
$ws.Range("M46").Value = "8/3/2018 19:48:00"
$ws.Range("M47").Value = "8/3/2018 19:49:00"
$ws.Range("M48").Value = "8/3/2018 19:49:00"
$ws.Range("M49").Value = "8/3/2018 19:49:00"
$ws.Range("M50").Value = "8/3/2018 19:49:00"
$ws.Range("M51").Value = "8/3/2018 19:49:00"
$ws.Range("M52").Value = "8/3/2018 19:50:00"
$ws.Range("M53").Value = "8/3/2018 19:50:00"
$ws.Range("M54").Value = "8/3/2018 19:50:00"
$ws.Range("M55").Value = "8/3/2018 19:50:00"
$ws.Range("M56").Value = "8/3/2018 19:50:00"
$ws.Range("M57").Value = "8/3/2018 19:51:00"
$ws.Range("M58").Value = "8/3/2018 19:51:00"
$ws.Range("M59").Value = "8/3/2018 19:51:00"
$ws.Range("M60").Value = "8/3/2018 19:51:00"
$ws.Range("M61").Value = "8/3/2018 19:51:00"
$ws.Range("M62").Value = "8/3/2018 19:51:00"
$ws.Range("M63").Value = "8/3/2018 19:51:00"
$ws.Range("M64").Value = "8/3/2018 19:51:00"
$ws.Range("M65").Value = "8/3/2018 19:51:00"
$ws.Range("M66").Value = "8/3/2018 19:52:00"
$ws.Range("M67").Value = "8/3/2018 19:52:00"
$ws.Range("M68").Value = "8/3/2018 19:52:00"
$ws.Range("M69").Value = "8/3/2018 19:53:00"
$ws.Range("M70").Value = "8/3/2018 19:53:00"
$ws.Range("M71").Value = "8/3/2018 19:53:00"
$ws.Range("M72").Value = "8/3/2018 19:55:00"
$ws.Range("M73").Value = "8/3/2018 19:55:00"
$ws.Range("M74").Value = "8/3/2018 19:56:00"
$ws.Range("M75").Value = "8/3/2018 19:56:00"
$ws.Range("M76").Value = "8/3/2018 19:57:00"
$ws.Range("M77").Value = "8/3/2018 19:57:00"
$ws.Range("M78").Value = "8/3/2018 19:57:00"
$ws.Range("M79").Value = "8/3/2018 19:57:00"
$ws.Range("M80").Value = "8/3/2018 19:57:00"
$ws.Range("M81").Value = "8/3/2018 19:57:00"
$ws.Range("M82").Value = "8/3/2018 19:57:00"
$ws.Range("M83").Value = "8/3/2018 19:57:00"
$ws.Range("M84").Value = "8/3/2018 19:57:00"
$ws.Range("M85").Value = "8/3/2018 19:58:00"
$ws.Range("M86").Value = "8/3/2018 19:59:00"
$ws.Range("M87").Value = "8/3/2018 19:59:00"
$ws.Range("M88").Value = "8/3/2018 20:00:00"
$ws.Range("M89").Value = "8/3/2018 20:01:00"
$ws.Range("M90").Value = "8/3/2018 20:01:00"
$ws.Range("M91").Value = "8/3/2018 20:01:00"
$ws.Range("M92").Value = "8/3/2018 20:02:00"
$ws.Range("M93").Value = "8/3/2018 20:02:00"
$ws.Range("M94").Value = "8/3/2018 20:02:00"
$ws.Range("M95").Value = "8/3/2018 20:02:00"
$ws.Range("M96").Value = "8/3/2018 20:02:00"
$ws.Range("M97").Value = "8/3/2018 20:02:00"
$ws.Range("M98").Value = "8/3/2018 20:02:00"
$ws.Range("M99").Value = "8/3/2018 20:03:00"
$ws.Range("M100").Value = "8/3/2018 20:03:00"
$ws.Range("M101").Value = "8/3/2018 20:03:00"
$ws.Range("M102").Value = "8/3/2018 20:03:00"
$ws.Range("M103").Value = "8/3/2018 20:03:00"
$ws.Range("M104").Value = "8/3/2018 20:03:00"
$ws.Range("M105").Value = "8/3/2018 20:03:00"
$ws.Range("M106").Value = "8/3/2018 20:04:00"
$ws.Range("M107").Value = "8/3/2018 20:04:00"
$ws.Range("M108").Value = "8/3/2018 20:04:00"
$ws.Range("M109").Value = "8/3/2018 20:05:00"
$ws.Range("M110").Value = "8/3/2018 20:05:00"
$ws.Range("M111").Value = "8/3/2018 20:05:00"
$ws.Range("M112").Value = "8/3/2018 20:05:00"
$ws.Range("M113").Value = "8/3/2018 20:05:00"
$ws.Range("M114").Value = "8/3/2018 20:05:00"
$ws.Range("M115").Value = "8/3/2018 20:06:00"
$ws.Range("M116").Value = "8/3/2018 20:07:00"
$ws.Range("M117").Value = "8/3/2018 20:08:00"
$ws.Range("M118").Value = "8/3/2018 20:08:00"
$ws.Range("M119").Value = "8/3/2018 20:08:00"
$ws.Range("M120").Value = "8/3/2018 20:08:00"
$ws.Range("M121").Value = "8/3/2018 20:10:00"
$ws.Range("M122").Value = "8/3/2018 20:10:00"
$ws.Range("M123").Value = "8/3/2018 20:10:00"
$ws.Range("M124").Value = "8/3/2018 20:10:00"
$ws.Range("M125").Value = "8/3/2018 20:11:00"
$ws.Range("M126").Value = "8/3/2018 20:11:00"
$ws.Range("M127").Value = "8/3/2018 20:11:00"
$ws.Range("M128").Value = "8/3/2018 20:11:00"
$ws.Range("M129").Value = "8/3/2018 20:11:00"
$ws.Range("M130").Value = "8/3/2018 20:11:00"
$ws.Range("M131").Value = "8/3/2018 20:11:00"
$ws.Range("M132").Value = "8/3/2018 20:11:00"
$ws.Range("M133").Value = "8/3/2018 20:12:00"
$ws.Range("M134").Value = "8/3/2018 20:12:00"
$ws.Range("M135").Value = "8/3/2018 20:12:00"
$ws.Range("M136").Value = "8/3/2018 20:12:00"
$ws.Range("M137").Value = "8/3/2018 20:13:00"
$ws.Range("M138").Value = "8/3/2018 20:13:00"
$ws.Range("M139").Value = "8/3/2018 20:13:00"
$ws.Range("M140").Value = "8/3/2018 20:13:00"
$ws.Range("M141").Value = "8/3/2018 20:13:00"
$ws.Range("M142").Value = "8/3/2018 20:13:00"
$ws.Range("M143").Value = "8/3/2018 20:13:00"
$ws.Range("M144").Value = "8/3/2018 20:14:00"
$ws.Range("M145").Value = "8/3/2018 20:15:00"
$ws.Range("M146").Value = "8/3/2018 20:32:00"
$ws.Range("M147").Value = "8/3/2018 20:32:00"
$ws.Range("M148").Value = "8/3/2018 20:32:00"
$ws.Range("M149").Value = "8/3/2018 20:32:00"
$ws.Range("M150").Value = "8/3/2018 20:32:00"
$ws.Range("M151").Value = "8/3/2018 20:32:00"
$ws.Range("M152").Value = "10/15/2018 14:43:00"
$ws.Range("M153").Value = "10/15/2018 14:43:00"
$ws.Range("M154").Value = "10/15/2018 14:43:00"
$ws.Range("M155").Value = "10/15/2018 14:44:00"
$ws.Range("M156").Value = "10/15/2018 14:44:00"
$ws.Range("M157").Value = "10/15/2018 14:44:00"
$ws.Range("M158").Value = "10/15/2018 14:46:00"
$ws.Range("M159").Value = "10/15/2018 14:47:00"
$ws.Range("M160").Value = "10/15/2018 14:49:00"
$ws.Range("M161").Value = "10/15/2018 14:49:00"
$ws.Range("M162").Value = "10/15/2018 14:52:00"
$ws.Range("M163").Value = "10/15/2018 14:52:00"
$ws.Range("M164").Value = "10/15/2018 14:52:00"
$ws.Range("M165").Value = "10/15/2018 14:52:00"
$ws.Range("M166").Value = "11/8/2018 11:38:00"
$ws.Range("M167").Value = "11/8/2018 11:38:00"
$ws.Range("M168").Value = "11/8/2018 11:38:00"
$ws.Range("M169").Value = "11/8/2018 11:38:00"
$ws.Range("M170").Value = "11/8/2018 11:39:00"
$ws.Range("M171").Value = "11/8/2018 11:39:00"
$ws.Range("M172").Value = "11/8/2018 11:39:00"
$ws.Range("M173").Value = "11/8/2018 11:39:00"
$ws.Range("M174").Value = "11/8/2018 11:39:00"
$ws.Range("M175").Value = "11/8/2018 11:39:00"
$ws.Range("M176").Value = "11/8/2018 11:40:00"
$ws.Range("M177").Value = "11/8/2018 11:40:00"
$ws.Range("M178").Value = "11/8/2018 11:40:00"
$ws.Range("M179").Value = "11/8/2018 14:29:00"
$ws.Range("M180").Value = "11/12/2018 13:33:00"
$ws.Range("M181").Value = "11/12/2018 13:33:00"
$ws.Range("M182").Value = "11/12/2018 13:33:00"
$ws.Range("M183").Value = "1/29/2019 16:47:09"
$ws.Range("M184").Value = "8/22/2019 14:19:16"
$ws.Range("M185").Value = "8/22/2019 14:19:20"
$ws.Range("M186").Value = "8/22/2019 14:19:48"

# Swap theme accent1 / accent5 colors
$theme = $wb.Theme
$scheme = $theme.ThemeColorScheme
$scheme.Colors(5).RGB  = 13998939   # accent1 -> 5B9BD5
$scheme.Colors(9).RGB  = 12874308   # accent5 -> 4472C4
